# Kilimanjaro Weekly Scoreboard - append the latest week's workout rows
# (Participant, Date, Workout Type, Total Duration, Total Distance,
#  Total Elevation, Zone 1-5, Workout Level, Week) for rows 127-132.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{A="Matt";    B=45472; C="Run";     D=51;  E=4.82; F=210; G=2;   H=35; I=9;  J=1; K=0; L="Agile Antelope";    M=3}
    @{A="Eric";    B=45472; C="Workout"; D=81;  E=0;    F=0;   G=9;   H=48; I=15; J=9; K=0; L="Sauntering Hippo";  M=3}
    @{A="Matt";    B=45472; C="Walk";    D=16;  E=0.59; F=10;  G=2;   H=0;  I=0;  J=0; K=0; L="Agile Antelope";    M=3}
    @{A="Steven";  B=45472; C="Walk";    D=140; E=5.57; F=856; G=140; H=0;  I=0;  J=0; K=0; L="Wily Hyena";        M=3}
    @{A="Phil";    B=45472; C="Run";     D=31;  E=3.12; F=203; G=0;   H=13; I=14; J=2; K=0; L="Sauntering Hippo";  M=3}
    @{A="Phil";    B=45472; C="Workout"; D=32;  E=0;    F=0;   G=13;  H=18; I=1;  J=0; K=0; L="Sauntering Hippo";  M=3}
)

# Existing dated cells (e.g. B126) use a short-date number format; grab it so
# the appended Date column (B) matches the rest of the table.
$ws.Range("B126").Copy()

$r = 127
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value  = $row.A
    $ws.Cells.Item($r, 2).Value  = $row.B
    $ws.Cells.Item($r, 3).Value  = $row.C
    $ws.Cells.Item($r, 4).Value  = $row.D
    $ws.Cells.Item($r, 5).Value  = $row.E
    $ws.Cells.Item($r, 6).Value  = $row.F
    $ws.Cells.Item($r, 7).Value  = $row.G
    $ws.Cells.Item($r, 8).Value  = $row.H
    $ws.Cells.Item($r, 9).Value  = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $r++
}

# Re-apply the date format to the newly added Date column cells.
$ws.Range("B127:B132").PasteSpecial(-4122)

# Move the active selection to the row right after the new data, like the
# workbook was left after scrolling down to review the appended entries.
$null = $ws.Range("A133").Select()
